$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.946.91'
$ws.Range('E2').Value = '  +5.58%  '
$ws.Range('D3').Value = '3.280.22'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '406.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('D7').Value = '3.275.69'
$ws.Range('E7').Value = '  +0.91%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.566'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.618'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.113'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +17.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '38.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').Value = '3.815.66'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '3.292.81'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').Value = '59.942.21'
$ws.Range('E18').Value = '  +5.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.982'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000113'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '294.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.42'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +16.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.39'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.31'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.74%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '133.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.292'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.119'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.04%  '
$ws.Range('D50').Value = '2.124.39'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('D51').Value = '3.616.73'
$ws.Range('E51').Value = '  +1.28%  '
